# modified stop_iv_chart: roll forward expiry dates by one week for the
# weekly-expiry rows (NIFTY weekly, BANKNIFTY, FINNIFTY, MIDCPNIFTY)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 45505
$ws.Range("B3").Value = 45512

$ws.Range("B6").Value = 45504
$ws.Range("B7").Value = 45511

$ws.Range("B8").Value = 45510
$ws.Range("B9").Value = 45517

$ws.Range("B10").Value = 45509
$ws.Range("B11").Value = 45516

$wb.Save()
